$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the affected rows
$ws.Range("F10").Value = 0
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = -2
$ws.Range("F21").Value = 0
$ws.Range("F24").Value = -5
$ws.Range("F26").Value = -9
$ws.Range("F28").Value = -1
$ws.Range("F32").Value = -7
$ws.Range("F33").Value = -3
$ws.Range("F38").Value = -2
$ws.Range("F40").Value = -2
